$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Rows.Item(5).Insert()
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "LP solver (linprog or gurobi)"
$ws.Range("B5").Value = "gurobi"
$ws.Range("A5").HorizontalAlignment = -4131
